$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 16
$ws.Range('A16').Value = 'V-1770523162330'
$ws.Range('B16').NumberFormat = '@'
$ws.Range('B16').Value = '7/2/2026'
$ws.Range('B16').Style = 'Normal'
$ws.Range('C16').Value = '10:59 p. m.'
$ws.Range('D16').Value = 'Stiven'
$ws.Range('E16').Value = 'Cerveza Corona (x1), Ron viejo de caldas (5años) botella (x1), Aguardiente Amarillo Media (x1), Aguardiente Amarillo Botella (x1)'
$ws.Range('F16').Value = 320000
$ws.Range('G16').Value = 0

# Row 17
$ws.Range('A17').Value = 'V-1770604719876'
$ws.Range('B17').NumberFormat = '@'
$ws.Range('B17').Value = '8/2/2026'
$ws.Range('B17').Style = 'Normal'
$ws.Range('C17').Value = '09:38 p. m.'
$ws.Range('D17').Value = 'Laura'
$ws.Range('E17').Value = 'Ron viejo de caldas (5años) botella (x1)'
$ws.Range('F17').Value = 132000
$ws.Range('G17').Value = 2

# Row 18
$ws.Range('A18').Value = 'V-1770604776176'
$ws.Range('B18').NumberFormat = '@'
$ws.Range('B18').Value = '8/2/2026'
$ws.Range('B18').Style = 'Normal'
$ws.Range('C18').Value = '09:39 p. m.'
$ws.Range('D18').Value = 'Laura'
$ws.Range('E18').Value = 'Cerveza Corona (x4)'
$ws.Range('F18').Value = 40000
$ws.Range('G18').Value = 0

# Row 19
$ws.Range('A19').Value = 'V-1770604832152'
$ws.Range('B19').NumberFormat = '@'
$ws.Range('B19').Value = '8/2/2026'
$ws.Range('B19').Style = 'Normal'
$ws.Range('C19').Value = '09:40 p. m.'
$ws.Range('D19').Value = 'Laura'
$ws.Range('E19').Value = 'Ron viejo de caldas (5años) botella (x1), Aguardiente Amarillo Botella (x1)'
$ws.Range('F19').Value = 240000
$ws.Range('G19').Value = 0

# Row 20
$ws.Range('A20').Value = 'V-1770604845907'
$ws.Range('B20').NumberFormat = '@'
$ws.Range('B20').Value = '8/2/2026'
$ws.Range('B20').Style = 'Normal'
$ws.Range('C20').Value = '09:40 p. m.'
$ws.Range('D20').Value = 'Laura'
$ws.Range('E20').Value = 'Cerveza Corona (x1)'
$ws.Range('F20').Value = 10000
$ws.Range('G20').Value = 0

# Row 21
$ws.Range('A21').Value = 'V-1770604879912'
$ws.Range('B21').NumberFormat = '@'
$ws.Range('B21').Value = '8/2/2026'
$ws.Range('B21').Style = 'Normal'
$ws.Range('C21').Value = '09:41 p. m.'
$ws.Range('D21').Value = 'Laura'
$ws.Range('E21').Value = 'Cerveza Corona (x3)'
$ws.Range('F21').Value = 30000
$ws.Range('G21').Value = 0

# Row 22
$ws.Range('A22').Value = 'V-1770604924749'
$ws.Range('B22').NumberFormat = '@'
$ws.Range('B22').Value = '8/2/2026'
$ws.Range('B22').Style = 'Normal'
$ws.Range('C22').Value = '09:42 p. m.'
$ws.Range('D22').Value = 'Laura'
$ws.Range('E22').Value = 'Ron viejo de caldas (5años) botella (x1), Cerveza Corona (x1)'
$ws.Range('F22').Value = 142000
$ws.Range('G22').Value = 0

# Row 23
$ws.Range('A23').Value = 'V-1770604965224'
$ws.Range('B23').NumberFormat = '@'
$ws.Range('B23').Value = '8/2/2026'
$ws.Range('B23').Style = 'Normal'
$ws.Range('C23').Value = '09:42 p. m.'
$ws.Range('D23').Value = 'Laura'
$ws.Range('E23').Value = 'Cerveza Corona (x10)'
$ws.Range('F23').Value = 100000
$ws.Range('G23').Value = 0

# Row 24
$ws.Range('A24').Value = 'V-1770605455118'
$ws.Range('B24').NumberFormat = '@'
$ws.Range('B24').Value = '8/2/2026'
$ws.Range('B24').Style = 'Normal'
$ws.Range('C24').Value = '09:50 p. m.'
$ws.Range('D24').Value = 'Laura'
$ws.Range('E24').Value = 'Aguardiente Amarillo Botella (x1)'
$ws.Range('F24').Value = 108000
$ws.Range('G24').Value = 0

# Row 25
$ws.Range('A25').Value = 'V-1770605530870'
$ws.Range('B25').NumberFormat = '@'
$ws.Range('B25').Value = '8/2/2026'
$ws.Range('B25').Style = 'Normal'
$ws.Range('C25').Value = '09:52 p. m.'
$ws.Range('D25').Value = 'Martha'
$ws.Range('E25').Value = 'Cerveza Corona (x1), Ron viejo de caldas (5años) botella (x1)'
$ws.Range('F25').Value = 142000
$ws.Range('G25').Value = 0

# Row 26
$ws.Range('A26').Value = 'V-1770763614972'
$ws.Range('B26').NumberFormat = '@'
$ws.Range('B26').Value = '10/2/2026'
$ws.Range('B26').Style = 'Normal'
$ws.Range('C26').Value = '05:46 p. m.'
$ws.Range('D26').Value = 'Laura'
$ws.Range('E26').Value = 'Ron viejo de caldas (5años) botella (x1)'
$ws.Range('F26').Value = 135000
$ws.Range('G26').Value = 2

# Row 27
$ws.Range('A27').Value = 'V-1770783380241'
$ws.Range('B27').NumberFormat = '@'
$ws.Range('B27').Value = '10/2/2026'
$ws.Range('B27').Style = 'Normal'
$ws.Range('C27').Value = '11:16 p. m.'
$ws.Range('D27').Value = 'Laura'
$ws.Range('E27').Value = 'Aguardiente Amarillo Media (x1)'
$ws.Range('F27').Value = 70000
$ws.Range('G27').Value = 2

# Row 28
$ws.Range('A28').Value = 'V-1770784241842'
$ws.Range('B28').NumberFormat = '@'
$ws.Range('B28').Value = '10/2/2026'
$ws.Range('B28').Style = 'Normal'
$ws.Range('C28').Value = '11:30 p. m.'
$ws.Range('D28').Value = 'Martha'
$ws.Range('E28').Value = 'Aguardiente Amarillo Caja (x1)'
$ws.Range('F28').Value = 122000
$ws.Range('G28').Value = 0
